# Split the sentence
#   "Do I want to use the return keyword to store a code as value in the function?"
# into
#   "Do I want to use the return keyword to store a code as the main code in the function?"
# and relocate the document's single "_GoBack" bookmark (Word keeps only one) from
# its old spot near "...xecuted its stored [BM] codes." to the newly-edited spot
# right before "in the function?" - exactly as real Word does after an edit.
#
# Note: this runtime normalises *all* runs of a paragraph into a single run
# whenever a Range's Text/InsertBefore/InsertAfter/Delete touches that paragraph,
# so run boundaries we want to keep are protected with temporary bookmarks
# (Bookmarks.Add splits a run in place without touching sibling runs) before any
# text-mutating call happens nearby.

$d = $word.ActiveDocument

$oldPhrase = "store a code as value in the function?"
$keepPrefix = "store a code as "        # stays in its own run (rsidR 00865EC2)
$newMiddle  = "the main code"           # brand new run
$tailText   = "in the function?"        # trailing run, gets the relocated bookmark before it

# --- locate the sentence fragment we need to touch -----------------------
$full = $d.Content.Text
$phraseStart = $full.IndexOf($oldPhrase)
if ($phraseStart -lt 0) { throw "Could not find target phrase" }

$splitPos = $phraseStart + $keepPrefix.Length

# --- 1) protect the "store a code as " / "value..." run boundary ---------
# (placed before the Find/Replace below so that replace's run-merge can't
# reach back across this point)
$protectRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("ZZ_TempProtect", $protectRange) | Out-Null

# --- 2) replace the single word "value" with "the main code" -------------
# Scope the Find to just this remainder of the sentence so no other
# occurrence of "value" elsewhere in the document is touched.
$full = $d.Content.Text
$remainderStart = $full.IndexOf("value " + $tailText)
$remainderEnd = $remainderStart + ("value " + $tailText).Length
$remainderRange = $d.Range($remainderStart, $remainderEnd)
$remainderRange.Find.Execute("value", $true, $false, $false, $false, $false, `
    $true, 0, $false, $newMiddle, 1) | Out-Null

# --- 3) split "the main code" / " " / "in the function?" apart -----------
$full = $d.Content.Text
$mergedStart = $full.IndexOf($newMiddle + " " + $tailText)
$afterMiddle = $mergedStart + $newMiddle.Length
$afterSpace = $afterMiddle + 1

$midSplitRange = $d.Range($afterMiddle, $afterMiddle)
$d.Bookmarks.Add("ZZ_TempSplit", $midSplitRange) | Out-Null

# --- 4) drop the (sole) "_GoBack" bookmark exactly where the edit ended --
# Adding it under this name automatically removes it from its previous
# location elsewhere in the document (Word allows only one "_GoBack").
$goBackRange = $d.Range($afterSpace, $afterSpace)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null

# --- 5) clean up the scaffolding bookmarks (run splits remain in place) --
$d.Bookmarks("ZZ_TempProtect").Delete()
$d.Bookmarks("ZZ_TempSplit").Delete()
